# Applies the tracked changes described by the commit diff:
#   1. Wrap "immigration:bullied" also in a grammar-check span
#      (adds <w:proofErr w:type="gramStart"/> / "gramEnd" around the
#      existing spell-check span) in the "Denmark is less integrative"
#      paragraph.
#   2. Append a brand-new paragraph right after that one with the
#      "A thing we want to specify..." text (red), including its own
#      gramStart/gramEnd proofErr markers.
#   3. Move the <w:lastRenderedPageBreak/> marker from the run that
#      starts "Immigration:hisced" to the run that starts "Escs_"
#      (earlier in the document).
#
# Because <w:proofErr/> and <w:lastRenderedPageBreak/> markers have no
# direct Word object-model surface (they are cosmetic artifacts of
# Word's proofing/pagination passes, not editable properties), each
# touched paragraph is rewritten wholesale with Range.InsertXML using
# the exact WordprocessingML for that paragraph - identical to the
# original except for the specific marker elements the diff adds,
# removes or relocates. Every other attribute/run is reproduced
# unchanged so the rest of the paragraph round-trips byte-for-byte.

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Get-ParagraphByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- 1 & 2: the Denmark/"immigration:bullied" paragraph -------------------
# Replace it with itself plus the new gramStart/gramEnd pair, then also
# insert the brand new paragraph right after it (still inside the same
# InsertXML call so both paragraphs land together, in one pass).

$denmarkPara = Get-ParagraphByText $d "Denmark is a less integrative country"
$denmarkXml = @"
<w:p $wns w14:paraId="13A18EEC" w14:textId="20324F73" w:rsidR="00800282" w:rsidRDefault="00800282" w:rsidP="002617B4"><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Furthermore, we observed that Denmark is a less integrative country: the interaction </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>immigration:bullied</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> has a negative effect on the scores, </w:t></w:r><w:r w:rsidR="008953F7"><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>underlining that an immigrant student who feels bullied at school is penalized.</w:t></w:r></w:p><w:p $wns><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">A thing we want to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>specify</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> and highlight is that all the models we built up are models with an R^2 that reaches the value of 0.30 to the maximum. Indeed, we know that the scores of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>a students</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> are not only due to the situations he lives in and things that doesn&#8217;t depend on him, but it is mostly made of his hard work and study. </w:t></w:r></w:p>
"@

$denmarkPara.Range.InsertXML($denmarkXml)

# --- 3: move <w:lastRenderedPageBreak/> from "Immigration:hisced" to "Escs_" ---

$escsPara = Get-ParagraphByText $d "Escs_status : positive effect (+28 circa)"
$escsXml = @"
<w:p $wns w14:paraId="2C9FEE10" w14:textId="1AFFF265" w:rsidR="00497760" w:rsidRPr="00927CF6" w:rsidRDefault="00497760" w:rsidP="004016E7"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:lastRenderedPageBreak/><w:t>Escs_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t>status</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> positive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t>effect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> (+28 circa) -&gt; se stai meglio economicamente e socialmente allora vai meglio </w:t></w:r></w:p>
"@

$escsPara.Range.InsertXML($escsXml)

$immPara = Get-ParagraphByText $d "Immigration:hisced : negative effect (-1,47)"
$immXml = @"
<w:p $wns w14:paraId="7131698F" w14:textId="359F957C" w:rsidR="008F2D0E" w:rsidRPr="00927CF6" w:rsidRDefault="00C51A44" w:rsidP="004016E7"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t>Immigration:hisced</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> : neg</w:t></w:r><w:r w:rsidR="00203AB1" w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t>a</w:t></w:r><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">tive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t>effect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00927CF6"><w:rPr><w:highlight w:val="green"/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> (-1,47) -&gt; per un immigrato i cui genitori sono poco istruiti risulta pi&#249; difficile andare bene a scuola </w:t></w:r></w:p>
"@

$immPara.Range.InsertXML($immXml)

Write-Host "Edits applied."
